$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '27.602.81'
$ws.Range("E2").Value = '  -0.44%  '
Set-TextValue $ws.Range("D3") '1.860.75'
$ws.Range("E3").Value = '  -0.87%  '
Set-TextValue $ws.Range("D4") '1.012'
$ws.Range("E4").Value = '  +0.74%  '
Set-TextValue $ws.Range("D5") '333.79'
$ws.Range("E5").Value = '  +0.34%  '
Set-TextValue $ws.Range("D6") '1.010'
$ws.Range("E6").Value = '  +0.66%  '
Set-TextValue $ws.Range("D7") '0.4669'
$ws.Range("E7").Value = '  -0.84%  '
Set-TextValue $ws.Range("D8") '0.3893'
$ws.Range("E8").Value = '  -1.59%  '
Set-TextValue $ws.Range("D9") '45.64'
$ws.Range("E9").Value = '  -4.56%  '
Set-TextValue $ws.Range("D10") '0.07967'
$ws.Range("E10").Value = '  -0.98%  '
Set-TextValue $ws.Range("D11") '1.001'
$ws.Range("E11").Value = '  -2.99%  '
Set-TextValue $ws.Range("D12") '21.62'
$ws.Range("E12").Value = '  -2.51%  '
Set-TextValue $ws.Range("D13") '1.855.45'
$ws.Range("E13").Value = '  -1.66%  '
Set-TextValue $ws.Range("D14") '5.974'
$ws.Range("E14").Value = '  -0.05%  '
Set-TextValue $ws.Range("D15") '7.218'
$ws.Range("E15").Value = '  +1.41%  '
Set-TextValue $ws.Range("D16") '1.012'
$ws.Range("E16").Value = '  +0.60%  '
Set-TextValue $ws.Range("D17") '87.90'
$ws.Range("E17").Value = '  +0.94%  '
Set-TextValue $ws.Range("D18") '0.06714'
$ws.Range("E18").Value = '  +0.59%  '
Set-TextValue $ws.Range("D19") '0.00001040'
$ws.Range("E19").Value = '  -0.79%  '
Set-TextValue $ws.Range("D20") '16.93'
$ws.Range("E20").Value = '  -1.71%  '
Set-TextValue $ws.Range("D21") '1.011'
$ws.Range("E21").Value = '  +0.71%  '
Set-TextValue $ws.Range("D22") '27.576.05'
$ws.Range("E22").Value = '  -0.60%  '
Set-TextValue $ws.Range("D23") '5.441'
$ws.Range("E23").Value = '  -1.52%  '
Set-TextValue $ws.Range("D24") '10.85'
$ws.Range("E24").Value = '  -1.34%  '
Set-TextValue $ws.Range("D25") '2.305'
$ws.Range("E25").Value = '  -0.18%  '
Set-TextValue $ws.Range("D26") '2.080.28'
$ws.Range("E26").Value = '  -1.27%  '
Set-TextValue $ws.Range("D27") '158.53'
$ws.Range("E27").Value = '  -0.51%  '
Set-TextValue $ws.Range("D28") '19.71'
$ws.Range("E28").Value = '  -2.20%  '
Set-TextValue $ws.Range("D29") '2.124'
$ws.Range("E29").Value = '  +1.02%  '
Set-TextValue $ws.Range("D30") '5.379'
$ws.Range("E30").Value = '  -3.49%  '
Set-TextValue $ws.Range("D31") '121.17'
$ws.Range("E31").Value = '  -0.47%  '
Set-TextValue $ws.Range("D32") '0.9698'
$ws.Range("E32").Value = '  -1.19%  '
Set-TextValue $ws.Range("D33") '0.09447'
$ws.Range("E33").Value = '  -0.72%  '
Set-TextValue $ws.Range("D34") '3.644'
$ws.Range("E34").Value = '  +1.26%  '
Set-TextValue $ws.Range("D35") '5.284'
$ws.Range("E35").Value = '  -1.23%  '
Set-TextValue $ws.Range("D36") '1.325'
$ws.Range("E36").Value = '  -8.30%  '
Set-TextValue $ws.Range("D37") '0.06029'
$ws.Range("E37").Value = '  -1.48%  '
Set-TextValue $ws.Range("D38") '0.02216'
$ws.Range("E38").Value = '  -2.03%  '
Set-TextValue $ws.Range("D39") '1.192'
$ws.Range("E39").Value = '  -3.21%  '
Set-TextValue $ws.Range("D40") '8.197'
$ws.Range("E40").Value = '  +0.71%  '
Set-TextValue $ws.Range("D41") '1.010'
$ws.Range("E41").Value = '  +0.69%  '
Set-TextValue $ws.Range("D42") '0.5905'
$ws.Range("E42").Value = '  -1.86%  '
Set-TextValue $ws.Range("D43") '0.1874'
$ws.Range("E43").Value = '  -1.38%  '
Set-TextValue $ws.Range("D44") '10.21'
$ws.Range("E44").Value = '  -0.64%  '
Set-TextValue $ws.Range("D45") '1.251'
$ws.Range("E45").Value = '  -1.18%  '
Set-TextValue $ws.Range("D46") '0.5613'
$ws.Range("E46").Value = '  -1.61%  '
Set-TextValue $ws.Range("D47") '12.09'
$ws.Range("E47").Value = '  -1.17%  '
Set-TextValue $ws.Range("D48") '1.913'
$ws.Range("E48").Value = '  -1.90%  '
Set-TextValue $ws.Range("D49") '3.269'
$ws.Range("E49").Value = '  -3.62%  '
Set-TextValue $ws.Range("D50") '0.06758'
$ws.Range("E50").Value = '  -2.33%  '
Set-TextValue $ws.Range("D51") '112.55'
$ws.Range("E51").Value = '  -1.31%  '
